$wb = $excel.ActiveWorkbook

# --- 1) Metadata sheet: bump the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- 2) Elements sheet: add the new "Mapping: Spécification métier..." column (AL) ---
$ws = $wb.Worksheets.Item("Elements")

# Header cell AL1 - copy formatting from AK1 (bold header style) then set text
$ws.Range("AK1").Copy()
$ws.Range("AL1").PasteSpecial(-4122)
$ws.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR FinancialHelpType"

# Data cells AL2:AL6 - copy formatting from AK2:AK6 (data row style)
$ws.Range("AK2:AK6").Copy()
$ws.Range("AL2:AL6").PasteSpecial(-4122)

# Only row 6 (Extension.value[x]) gets a mapping value
$ws.Range("AL6").Value = "aideFinanciere"

# Match the column width used for the new column in the source workbook (~74.07 chars)
$ws.Columns.Item(38).ColumnWidth = 73.15
